# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list (column E, rows 16-22) is reordered from descending
# (2306..2212) to ascending (2212..2306), each period keeping its own
# "Valor Mora" (column F) value, and every "Salario Basico" (column G)
# value is updated from 1500000 to 1200000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (Periodo Mora, Valor Mora, Salario Basico)
$rows = @(
    @{ Row = 16; Periodo = "2212"; Mora = 60000; Salario = 1200000 },
    @{ Row = 17; Periodo = "2301"; Mora = 60000; Salario = 1200000 },
    @{ Row = 18; Periodo = "2302"; Mora = 60000; Salario = 1200000 },
    @{ Row = 19; Periodo = "2303"; Mora = 60000; Salario = 1200000 },
    @{ Row = 20; Periodo = "2304"; Mora = 60000; Salario = 1200000 },
    @{ Row = 21; Periodo = "2305"; Mora = 60000; Salario = 1200000 },
    @{ Row = 22; Periodo = "2306"; Mora = 35200; Salario = 1200000 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 5).Value = $item.Periodo   # column E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $item.Mora       # column F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $item.Salario    # column G - Salario Basico
}
